$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("B2").Value = 'Bitcoin'
$ws.Range("C2").Value = 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'
$ws.Range("D2").Value = '58.492.01'
$ws.Range("E2").Value = '  -4.00%  '

$ws.Range("B3").Value = 'Ethereum'
$ws.Range("C3").Value = 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'
$ws.Range("D3").Value = '2.550.98'
$ws.Range("E3").Value = '  -1.89%  '

$ws.Range("B4").Value = 'TetherUSD'
$ws.Range("C4").Value = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("B5").Value = 'BNB'
$ws.Range("C5").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D5").Value = '505.35'
$ws.Range("E5").Value = '  -3.54%  '

$ws.Range("B6").Value = 'Solana'
$ws.Range("C6").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D6").Value = '143.79'
$ws.Range("E6").Value = '  -7.23%  '

$ws.Range("B7").Value = 'USDC'
$ws.Range("C7").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  +0.12%  '

$ws.Range("B8").Value = 'XRP'
$ws.Range("C8").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D8").Value = '0.552'
$ws.Range("E8").Value = '  -7.31%  '

$ws.Range("B9").Value = 'LidoStakedEther'
$ws.Range("C9").Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range("D9").Value = '2.554.97'
$ws.Range("E9").Value = '  -2.00%  '

$ws.Range("B10").Value = 'Toncoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D10").Value = '6.18'
$ws.Range("E10").Value = '  -7.64%  '

$ws.Range("B11").Value = 'Dogecoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D11").Value = '0.101'
$ws.Range("E11").Value = '  -3.90%  '

$ws.Range("B12").Value = 'Cardano'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D12").Value = '0.331'
$ws.Range("E12").Value = '  -4.69%  '

$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D13").Value = '0.128'
$ws.Range("E13").Value = '  -0.87%  '

$ws.Range("B14").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C14").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D14").Value = '2.998.72'
$ws.Range("E14").Value = '  -1.75%  '

$ws.Range("B15").Value = 'WrappedBTC'
$ws.Range("C15").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D15").Value = '58.488.34'
$ws.Range("E15").Value = '  -4.04%  '

$ws.Range("B16").Value = 'Avalanche'
$ws.Range("C16").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D16").Value = '20.48'
$ws.Range("E16").Value = '  -5.72%  '

$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").Value = '0.0000133'
$ws.Range("E17").Value = '  -5.76%  '

$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").Value = '2.554.96'
$ws.Range("E18").Value = '  -1.75%  '

$ws.Range("B19").Value = 'Polkadot'
$ws.Range("C19").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D19").Value = '4.50'
$ws.Range("E19").Value = '  -5.48%  '

$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D20").Value = '333.16'
$ws.Range("E20").Value = '  -6.10%  '

$ws.Range("B21").Value = 'Chainlink'
$ws.Range("C21").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D21").Value = '10.04'
$ws.Range("E21").Value = '  -5.08%  '

$ws.Range("B22").Value = 'Dai'
$ws.Range("C22").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D22").Value = '0.996'
$ws.Range("E22").Value = '  -0.22%  '

$ws.Range("B23").Value = 'Uniswap'
$ws.Range("C23").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D23").Value = '5.93'
$ws.Range("E23").Value = '  -4.71%  '

$ws.Range("B24").Value = 'Litecoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D24").Value = '59.37'
$ws.Range("E24").Value = '  -2.54%  '

$ws.Range("B25").Value = 'Polygon'
$ws.Range("C25").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D25").Value = '0.406'
$ws.Range("E25").Value = '  -4.93%  '

$ws.Range("B26").Value = 'Binance-PegBSC-USD'
$ws.Range("C26").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D26").Value = '1.00'
$ws.Range("E26").Value = '  +0.51%  '

$ws.Range("B27").Value = 'Kaspa'
$ws.Range("C27").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D27").Value = '0.155'
$ws.Range("E27").Value = '  -6.92%  '

$ws.Range("B28").Value = 'PEPE'
$ws.Range("C28").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D28").Value = '0.0₃0774'
$ws.Range("E28").Value = '  -8.69%  '

$ws.Range("B29").Value = 'InternetComputer(DFINITY)'
$ws.Range("C29").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D29").Value = '6.85'
$ws.Range("E29").Value = '  -7.37%  '

$ws.Range("B30").Value = 'USDe'
$ws.Range("C30").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D30").Value = '1.00'
$ws.Range("E30").Value = '  -0.02%  '

$ws.Range("B31").Value = 'Aptos'
$ws.Range("C31").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D31").Value = '5.85'
$ws.Range("E31").Value = '  -7.26%  '

$ws.Range("B32").Value = 'EthereumClassic'
$ws.Range("C32").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D32").Value = '18.53'
$ws.Range("E32").Value = '  -4.44%  '

$ws.Range("B33").Value = 'Monero'
$ws.Range("C33").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D33").Value = '148.80'
$ws.Range("E33").Value = '  +0.68%  '

$ws.Range("B34").Value = 'PancakeSwap'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D34").Value = '1.54'
$ws.Range("E34").Value = '  -4.18%  '

$ws.Range("B35").Value = 'SuiNetwork'
$ws.Range("C35").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D35").Value = '0.909'
$ws.Range("E35").Value = '  -2.85%  '

$ws.Range("B36").Value = 'NEARProtocol'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D36").Value = '3.88'
$ws.Range("E36").Value = '  -7.27%  '

$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D37").Value = '1.10'
$ws.Range("E37").Value = '  -8.20%  '

$ws.Range("B38").Value = 'OKB'
$ws.Range("C38").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D38").Value = '35.85'
$ws.Range("E38").Value = '  -1.67%  '

$ws.Range("B39").Value = 'Fetch.AI'
$ws.Range("C39").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D39").Value = '0.816'
$ws.Range("E39").Value = '  -5.51%  '

$ws.Range("B40").Value = 'Bittensor'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D40").Value = '282.98'
$ws.Range("E40").Value = '  -1.82%  '

$ws.Range("B41").Value = 'Filecoin'
$ws.Range("C41").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D41").Value = '3.52'
$ws.Range("E41").Value = '  -7.49%  '

$ws.Range("B42").Value = 'Stacks'
$ws.Range("C42").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D42").Value = '1.37'
$ws.Range("E42").Value = '  -9.13%  '

$ws.Range("B43").Value = 'FirstDigitalUSD'
$ws.Range("C43").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D43").Value = '0.999'
$ws.Range("E43").Value = '  +0.23%  '

$ws.Range("B44").Value = 'Stellar'
$ws.Range("C44").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D44").Value = '0.0980'
$ws.Range("E44").Value = '  -3.61%  '

$ws.Range("B45").Value = 'Mantle'
$ws.Range("C45").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D45").Value = '0.606'
$ws.Range("E45").Value = '  -2.74%  '

$ws.Range("B46").Value = 'Hedera'
$ws.Range("C46").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D46").Value = '0.0532'
$ws.Range("E46").Value = '  -5.28%  '

$ws.Range("B47").Value = 'WhiteBITCoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D47").Value = '10.33'
$ws.Range("E47").Value = '  -0.01%  '

$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").Value = '18.57'
$ws.Range("E48").Value = '  -5.45%  '

$ws.Range("B49").Value = 'VeChain'
$ws.Range("C49").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D49").Value = '0.0226'
$ws.Range("E49").Value = '  -4.96%  '

$ws.Range("B50").Value = 'RenderToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D50").Value = '4.50'
$ws.Range("E50").Value = '  -9.59%  '

$ws.Range("B51").Value = 'Maker'
$ws.Range("C51").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D51").Value = '1.902.39'
$ws.Range("E51").Value = '  -3.19%  '
